$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.26, -0.073, -0.133, 37, 19, 97),
    @(0.114, 0.172, 0.507, -89, -2, -164),
    @(0.172, -0.077, 0.156, 88, 86, -172),
    @(-0.059, -0.222, 0.08500000000000001, 108, -10, -140),
    @(0.191, -0.379, 0.238, 16, 59, 55),
    @(-0.196, -0.002, 0.243, -51, -49, -13),
    @(-0.093, 0.24, 0.472, 92, -10, 75),
    @(-0.2, -0.09, 0.161, 115, -36, 136),
    @(0.116, -0.017, 0.148, -78, -7, -127),
    @(0.025, -0.135, 0.158, -95, 13, -149)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$wb.Save()
